# SimpleTriangleCppWinRT_UWP Readme.docx
#
# 1. Update the "compatible SDK" sentence from the Fall Creators Update
#    SDK (16299) to the April 2018 Update SDK (17134), merging the two
#    runs that made up the sentence into a single run.
# 2. That paragraph becomes a (visually unchanged) Heading1-styled
#    paragraph with spacing-before reset to 0 and direct character
#    formatting overrides that keep its original italic / 10pt look.
# 3. The stray "_GoBack" bookmark that Word had left on the "Using the
#    sample" heading moves to the blank paragraph right after the SDK
#    sentence.

$d = $word.ActiveDocument

# --- 1. Swap the SDK sentence -------------------------------------------
$d.Content.Find.Execute(
    "This sample is compatible with the Windows 10 Fall Creators Update SDK (16299)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This sample is compatible with the Windows 10 April 2018 Update SDK (17134)",
    2) | Out-Null

# --- 2. Re-style that paragraph as Heading1, preserving its look -------
$p = $d.Paragraphs.Item(2)
$p.Style = "Heading 1"
$p.Format.SpaceBefore = 0

$p.Range.Font.Italic = $true
$p.Range.Font.Color = -16777216    # wdColorAutomatic -> w:color w:val="auto"
$p.Range.Font.Size = 10            # w:sz (half-points) = 20
$p.Range.Font.SizeBi = 11          # w:szCs (half-points) = 22
$p.Range.Font.NameBi = "Times New Roman"

# --- 3. Move the "_GoBack" bookmark -------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$blank = $d.Paragraphs.Item(3)
$d.Bookmarks.Add("_GoBack", $blank.Range)
